$wb = $excel.ActiveWorkbook

# --- ccpd_post_st (sheet2 / rId2) ---
$ws = $wb.Worksheets.Item("ccpd_post_st")
$ws.Range("C2").Value = 80
$ws.Activate()
$ws.Range("C3").Select()

# --- uc_running_cfg (sheet3 / rId3) ---
$ws = $wb.Worksheets.Item("uc_running_cfg")
$ws.Range("C2").Value = 10
$ws.Activate()
$ws.Range("C2").Select()

# --- heater_cfg (sheet4 / rId4) ---
$ws = $wb.Worksheets.Item("heater_cfg")
$ws.Range("C2").Value = 20
$ws.Activate()
$ws.Range("C2").Select()

# --- bms_cfg (sheet5 / rId5) ---
$ws = $wb.Worksheets.Item("bms_cfg")
$ws.Range("C2").Value = 30
$ws.Activate()
$ws.Range("C2").Select()

# --- mainpump_cfg (sheet6 / rId6) ---
$ws = $wb.Worksheets.Item("mainpump_cfg")
$ws.Range("C2").Value = 40
$ws.Activate()
$ws.Range("C2").Select()

# --- valve_cfg (sheet7 / rId7) ---
$ws = $wb.Worksheets.Item("valve_cfg")
$ws.Range("C2").Value = 50
$ws.Activate()
$ws.Range("C2").Select()

# --- sensors_cfg (sheet8 / rId8) ---
$ws = $wb.Worksheets.Item("sensors_cfg")
$ws.Range("C2").Value = 60
$ws.Activate()
$ws.Range("C2").Select()

# --- pressure_cfg (sheet9 / rId9) ---
$ws = $wb.Worksheets.Item("pressure_cfg")
$ws.Range("C2").Value = 70
$ws.Activate()
$ws.Range("C2").Select()

# --- config (sheet1 / rId1) : ends up as the active tab ---
$ws = $wb.Worksheets.Item("config")
$ws.Activate()
$ws.Range("E7").Select()
